# Append the new rows (234-238) covering 22-26 April 2021, as per the
# "aggiornamento fino al 26/03" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44308, 10, 40, 368.4937816674343),
    @(44309, 3, 34, 313.2197144173192),
    @(44310, 3, 31, 285.5826807922617),
    @(44311, 0, 30, 276.3703362505758),
    @(44312, 2, 21, 193.459235375403)
)

$startRow = 234
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Mirror the style (date format + border + alignment) of the last
    # existing data row (233) onto column A of the new row, then set the
    # real values for all four columns.
    $ws.Range("A233").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
